$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value, and whether the value needs to be
# forced to Text (to stop Excel from re-interpreting numeric-looking strings,
# e.g. "573.90", as actual numbers and stripping the trailing zero / precision).
$updates = @(
    @{ Cell = "D2"; Value = "66.274.38"; ForceText = $false },
    @{ Cell = "E2"; Value = "  -0.40%  "; ForceText = $false },
    @{ Cell = "D3"; Value = "3.522.61"; ForceText = $false },
    @{ Cell = "E3"; Value = "  +1.48%  "; ForceText = $false },
    @{ Cell = "E4"; Value = "  +0.02%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "573.90"; ForceText = $true },
    @{ Cell = "E5"; Value = "  +5.41%  "; ForceText = $false },
    @{ Cell = "D6"; Value = "179.78"; ForceText = $true },
    @{ Cell = "E6"; Value = "  -4.06%  "; ForceText = $false },
    @{ Cell = "E7"; Value = "  +5.81%  "; ForceText = $false },
    @{ Cell = "E8"; Value = "  -0.01%  "; ForceText = $false },
    @{ Cell = "D9"; Value = "0.639"; ForceText = $true },
    @{ Cell = "E9"; Value = "  +1.87%  "; ForceText = $false },
    @{ Cell = "D10"; Value = "0.157"; ForceText = $true },
    @{ Cell = "E10"; Value = "  +5.21%  "; ForceText = $false },
    @{ Cell = "D11"; Value = "55.97"; ForceText = $true },
    @{ Cell = "E11"; Value = "  +3.06%  "; ForceText = $false },
    @{ Cell = "D12"; Value = "0.0000276"; ForceText = $true },
    @{ Cell = "E12"; Value = "  +4.10%  "; ForceText = $false },
    @{ Cell = "D13"; Value = "9.36"; ForceText = $true },
    @{ Cell = "E13"; Value = "  +0.77%  "; ForceText = $false },
    @{ Cell = "D14"; Value = "4.080.81"; ForceText = $false },
    @{ Cell = "E14"; Value = "  +1.51%  "; ForceText = $false },
    @{ Cell = "D15"; Value = "3.518.48"; ForceText = $false },
    @{ Cell = "E15"; Value = "  +1.46%  "; ForceText = $false },
    @{ Cell = "E16"; Value = "  +0.44%  "; ForceText = $false },
    @{ Cell = "D17"; Value = "18.44"; ForceText = $true },
    @{ Cell = "E17"; Value = "  +2.57%  "; ForceText = $false },
    @{ Cell = "D18"; Value = "66.228.45"; ForceText = $false },
    @{ Cell = "E18"; Value = "  -0.45%  "; ForceText = $false },
    @{ Cell = "D19"; Value = "12.08"; ForceText = $true },
    @{ Cell = "E19"; Value = "  +3.43%  "; ForceText = $false },
    @{ Cell = "E20"; Value = "  +2.52%  "; ForceText = $false },
    @{ Cell = "D21"; Value = "416.54"; ForceText = $true },
    @{ Cell = "E21"; Value = "  -0.74%  "; ForceText = $false },
    @{ Cell = "D22"; Value = "4.20"; ForceText = $true },
    @{ Cell = "E22"; Value = "  +8.65%  "; ForceText = $false },
    @{ Cell = "B23"; Value = "Litecoin"; ForceText = $false },
    @{ Cell = "C23"; Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; ForceText = $false },
    @{ Cell = "D23"; Value = "85.75"; ForceText = $true },
    @{ Cell = "E23"; Value = "  +2.00%  "; ForceText = $false },
    @{ Cell = "B24"; Value = "Toncoin"; ForceText = $false },
    @{ Cell = "C24"; Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; ForceText = $false },
    @{ Cell = "D24"; Value = "4.26"; ForceText = $true },
    @{ Cell = "E24"; Value = "  +2.76%  "; ForceText = $false },
    @{ Cell = "D25"; Value = "13.02"; ForceText = $true },
    @{ Cell = "E25"; Value = "  +10.28%  "; ForceText = $false },
    @{ Cell = "D26"; Value = "11.06"; ForceText = $true },
    @{ Cell = "E26"; Value = "  -0.07%  "; ForceText = $false },
    @{ Cell = "D27"; Value = "2.87"; ForceText = $true },
    @{ Cell = "E27"; Value = "  -0.08%  "; ForceText = $false },
    @{ Cell = "D28"; Value = "9.13"; ForceText = $true },
    @{ Cell = "E28"; Value = "  +4.49%  "; ForceText = $false },
    @{ Cell = "D29"; Value = "30.58"; ForceText = $true },
    @{ Cell = "E29"; Value = "  +2.34%  "; ForceText = $false },
    @{ Cell = "D30"; Value = "649.33"; ForceText = $true },
    @{ Cell = "E30"; Value = "  -0.34%  "; ForceText = $false },
    @{ Cell = "D31"; Value = "6.60"; ForceText = $true },
    @{ Cell = "E31"; Value = "  +0.12%  "; ForceText = $false },
    @{ Cell = "D32"; Value = "11.74"; ForceText = $true },
    @{ Cell = "E32"; Value = "  +0.97%  "; ForceText = $false },
    @{ Cell = "E33"; Value = "  +2.13%  "; ForceText = $false },
    @{ Cell = "D34"; Value = "0.157"; ForceText = $true },
    @{ Cell = "E34"; Value = "  +14.48%  "; ForceText = $false },
    @{ Cell = "D35"; Value = "59.60"; ForceText = $true },
    @{ Cell = "E35"; Value = "  +1.02%  "; ForceText = $false },
    @{ Cell = "E36"; Value = "  +0.45%  "; ForceText = $false },
    @{ Cell = "D37"; Value = "0.0₃0803"; ForceText = $false },
    @{ Cell = "E37"; Value = "  -0.60%  "; ForceText = $false },
    @{ Cell = "D38"; Value = "37.48"; ForceText = $true },
    @{ Cell = "E38"; Value = "  -1.53%  "; ForceText = $false },
    @{ Cell = "D39"; Value = "0.383"; ForceText = $true },
    @{ Cell = "E39"; Value = "  -0.94%  "; ForceText = $false },
    @{ Cell = "D40"; Value = "3.267.47"; ForceText = $false },
    @{ Cell = "E40"; Value = "  +9.27%  "; ForceText = $false },
    @{ Cell = "D41"; Value = "3.47"; ForceText = $true },
    @{ Cell = "E41"; Value = "  +4.27%  "; ForceText = $false },
    @{ Cell = "E42"; Value = "  -0.03%  "; ForceText = $false },
    @{ Cell = "D43"; Value = "2.92"; ForceText = $true },
    @{ Cell = "E43"; Value = "  +2.00%  "; ForceText = $false },
    @{ Cell = "D44"; Value = "3.35"; ForceText = $true },
    @{ Cell = "E44"; Value = "  -2.09%  "; ForceText = $false },
    @{ Cell = "D45"; Value = "2.54"; ForceText = $true },
    @{ Cell = "E45"; Value = "  -2.57%  "; ForceText = $false },
    @{ Cell = "D46"; Value = "0.0420"; ForceText = $true },
    @{ Cell = "E46"; Value = "  +1.70%  "; ForceText = $false },
    @{ Cell = "D47"; Value = "2.72"; ForceText = $true },
    @{ Cell = "E47"; Value = "  +1.26%  "; ForceText = $false },
    @{ Cell = "E48"; Value = "  +3.06%  "; ForceText = $false },
    @{ Cell = "D49"; Value = "8.75"; ForceText = $true },
    @{ Cell = "E49"; Value = "  +0.20%  "; ForceText = $false },
    @{ Cell = "D50"; Value = "139.39"; ForceText = $true },
    @{ Cell = "E50"; Value = "  +0.44%  "; ForceText = $false },
    @{ Cell = "D51"; Value = "2.47"; ForceText = $true },
    @{ Cell = "E51"; Value = "  +1.61%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Prefix with an apostrophe so Excel stores the value as literal text
        # instead of parsing it into a Double (which would round/alter the digits).
        $range.Value = "'" + $u.Value
        # Re-apply the plain "Normal" cell style so the quote-prefix text style
        # that Excel auto-assigns does not change the cells visible formatting.
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
